# Edit: split the "Converted to hex" run into "Converted " + "to decimal (from hex)"
# i.e. paragraph text goes from
#   "Converted to hex: [02, 01, 166, 86, 0]"
# to
#   "Converted to decimal (from hex): [02, 01, 166, 86, 0]"
# keeping the leading "Converted " / remaining colon-run formatting (italic on the
# "Converted "/"to decimal (from hex)" runs, non-italic on the ": " run) intact.

$p = $ppt.ActivePresentation

$targetSlide = $null
$targetShape = $null
$targetParaIndex = 0
$oldLead = "Converted to hex"
$searchNeedle = "Converted to hex"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shp = $slide.Shapes.Item($shi)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            $paraCount = $tr.Paragraphs().Count
            for ($pi = 1; $pi -le $paraCount; $pi++) {
                $para = $tr.Paragraphs($pi, 1)
                if ($para.Text.Contains($searchNeedle)) {
                    $targetSlide = $slide
                    $targetShape = $shp
                    $targetParaIndex = $pi
                }
            }
        }
    }
}

if ($targetParaIndex -eq 0) {
    Write-Host "Paragraph containing '$searchNeedle' not found"
} else {
    $tr = $targetShape.TextFrame.TextRange
    $para = $tr.Paragraphs($targetParaIndex, 1)

    # "Converted to hex" = "Converted " (10 chars) + "to hex" (6 chars), starting at
    # character 1 of the paragraph.
    $leadLen = $oldLead.Length          # 16 -> "Converted to hex"
    $keepLen = 10                        # "Converted "
    $restStart = $keepLen + 1            # 11
    $restLen = $leadLen - $keepLen       # 6 -> "to hex"

    # Sanity check before mutating.
    $check = $para.Characters(1, $leadLen).Text
    if ($check -ne $oldLead) {
        Write-Host "Unexpected text at target paragraph: [$check]"
    } else {
        $tail = $para.Characters($restStart, $restLen)
        $tail.Text = "to decimal (from hex)"
        Write-Host "Updated paragraph text: $($para.Text)"
    }
}
